# Natmi following Dr Hou advice
# Update ligand/receptor-expressing cell counts (E, K) and all dependent
# expression / specificity statistics (G, H, I, J, M, N, O, P, Q, R, S, T)
# on Sheet1 to reflect the revised cell-calling threshold.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 15.33750333333333
$ws.Range("H2").Value = 46.01251
$ws.Range("I2").Value = 0.1440483515229198
$ws.Range("J2").Value = 0.1440483515229198
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 1.776179
$ws.Range("N2").Value = 5.328537
$ws.Range("O2").Value = 0.009213114886297067
$ws.Range("P2").Value = 0.009213114886297067
$ws.Range("Q2").Value = 27.24215133309666
$ws.Range("R2").Value = 245.17936199787
$ws.Range("S2").Value = 0.001327134011762365
$ws.Range("T2").Value = 0.001327134011762365

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 15.33750333333333
$ws.Range("H3").Value = 46.01251
$ws.Range("I3").Value = 0.1440483515229198
$ws.Range("J3").Value = 0.1440483515229198
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 103.273595
$ws.Range("N3").Value = 309.820785
$ws.Range("O3").Value = 0.5356844639284184
$ws.Range("P3").Value = 0.5356844639284185
$ws.Range("Q3").Value = 1583.959107557817
$ws.Range("R3").Value = 14255.63196802035
$ws.Range("S3").Value = 0.07716446396532768
$ws.Range("T3").Value = 0.0771644639653277

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 15.33750333333333
$ws.Range("H4").Value = 46.01251
$ws.Range("I4").Value = 0.1440483515229198
$ws.Range("J4").Value = 0.1440483515229198
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 53.963124
$ws.Range("N4").Value = 161.889372
$ws.Range("O4").Value = 0.2799089849815219
$ws.Range("P4").Value = 0.2799089849815219
$ws.Range("Q4").Value = 827.65959422708
$ws.Range("R4").Value = 7448.93634804372
$ws.Range("S4").Value = 0.04032042786304196
$ws.Range("T4").Value = 0.04032042786304196

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 15.33750333333333
$ws.Range("H5").Value = 46.01251
$ws.Range("I5").Value = 0.1440483515229198
$ws.Range("J5").Value = 0.1440483515229198
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 33.77521133333333
$ws.Range("N5").Value = 101.325634
$ws.Range("O5").Value = 0.1751934362037625
$ws.Range("P5").Value = 0.1751934362037625
$ws.Range("Q5").Value = 518.0274164090378
$ws.Range("R5").Value = 4662.24674768134
$ws.Range("S5").Value = 0.02523632568278781
$ws.Range("T5").Value = 0.02523632568278782

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 64.92210766666668
$ws.Range("H6").Value = 194.766323
$ws.Range("I6").Value = 0.609742171429695
$ws.Range("J6").Value = 0.6097421714296949
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 1.776179
$ws.Range("N6").Value = 5.328537
$ws.Range("O6").Value = 0.009213114886297067
$ws.Range("P6").Value = 0.009213114886297067
$ws.Range("Q6").Value = 115.3132842732723
$ws.Range("R6").Value = 1037.819558459451
$ws.Range("S6").Value = 0.005617624676402021
$ws.Range("T6").Value = 0.00561762467640202

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 64.92210766666668
$ws.Range("H7").Value = 194.766323
$ws.Range("I7").Value = 0.609742171429695
$ws.Range("J7").Value = 0.6097421714296949
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 103.273595
$ws.Range("N7").Value = 309.820785
$ws.Range("O7").Value = 0.5356844639284184
$ws.Range("P7").Value = 0.5356844639284185
$ws.Range("Q7").Value = 6704.73945371373
$ws.Range("R7").Value = 60342.65508342356
$ws.Range("S7").Value = 0.326629408236866
$ws.Range("T7").Value = 0.326629408236866

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 64.92210766666668
$ws.Range("H8").Value = 194.766323
$ws.Range("I8").Value = 0.609742171429695
$ws.Range("J8").Value = 0.6097421714296949
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 53.963124
$ws.Range("N8").Value = 161.889372
$ws.Range("O8").Value = 0.2799089849815219
$ws.Range("P8").Value = 0.2799089849815219
$ws.Range("Q8").Value = 3503.399746357685
$ws.Range("R8").Value = 31530.59771721916
$ws.Range("S8").Value = 0.170672312305315
$ws.Range("T8").Value = 0.170672312305315

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 64.92210766666668
$ws.Range("H9").Value = 194.766323
$ws.Range("I9").Value = 0.609742171429695
$ws.Range("J9").Value = 0.6097421714296949
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 33.77521133333333
$ws.Range("N9").Value = 101.325634
$ws.Range("O9").Value = 0.1751934362037625
$ws.Range("P9").Value = 0.1751934362037625
$ws.Range("Q9").Value = 2192.757906647087
$ws.Range("R9").Value = 19734.82115982378
$ws.Range("S9").Value = 0.1068228262111119
$ws.Range("T9").Value = 0.1068228262111119

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 10.67805633333333
$ws.Range("H10").Value = 32.034169
$ws.Range("I10").Value = 0.1002872748488753
$ws.Range("J10").Value = 0.1002872748488753
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 1.776179
$ws.Range("N10").Value = 5.328537
$ws.Range("O10").Value = 0.009213114886297067
$ws.Range("P10").Value = 0.009213114886297067
$ws.Range("Q10").Value = 18.96613942008366
$ws.Range("R10").Value = 170.695254780753
$ws.Range("S10").Value = 0.000923958184816338
$ws.Range("T10").Value = 0.000923958184816338

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 10.67805633333333
$ws.Range("H11").Value = 32.034169
$ws.Range("I11").Value = 0.1002872748488753
$ws.Range("J11").Value = 0.1002872748488753
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 103.273595
$ws.Range("N11").Value = 309.820785
$ws.Range("O11").Value = 0.5356844639284184
$ws.Range("P11").Value = 0.5356844639284185
$ws.Range("Q11").Value = 1102.761265155852
$ws.Range("R11").Value = 9924.851386402665
$ws.Range("S11").Value = 0.0537223350662617
$ws.Range("T11").Value = 0.05372233506626171

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 10.67805633333333
$ws.Range("H12").Value = 32.034169
$ws.Range("I12").Value = 0.1002872748488753
$ws.Range("J12").Value = 0.1002872748488753
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 53.963124
$ws.Range("N12").Value = 161.889372
$ws.Range("O12").Value = 0.2799089849815219
$ws.Range("P12").Value = 0.2799089849815219
$ws.Range("Q12").Value = 576.2212779946519
$ws.Range("R12").Value = 5185.991501951868
$ws.Range("S12").Value = 0.02807130930951158
$ws.Range("T12").Value = 0.02807130930951158

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 10.67805633333333
$ws.Range("H13").Value = 32.034169
$ws.Range("I13").Value = 0.1002872748488753
$ws.Range("J13").Value = 0.1002872748488753
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 33.77521133333333
$ws.Range("N13").Value = 101.325634
$ws.Range("O13").Value = 0.1751934362037625
$ws.Range("P13").Value = 0.1751934362037625
$ws.Range("Q13").Value = 360.6536092875717
$ws.Range("R13").Value = 3245.882483588146
$ws.Range("S13").Value = 0.01756967228828563
$ws.Range("T13").Value = 0.01756967228828563

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 15.537021
$ws.Range("H14").Value = 46.611063
$ws.Range("I14").Value = 0.14592220219851
$ws.Range("J14").Value = 0.14592220219851
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 1.776179
$ws.Range("N14").Value = 5.328537
$ws.Range("O14").Value = 0.009213114886297067
$ws.Range("P14").Value = 0.009213114886297067
$ws.Range("Q14").Value = 27.596530422759
$ws.Range("R14").Value = 248.368773804831
$ws.Range("S14").Value = 0.001344398013316343
$ws.Range("T14").Value = 0.001344398013316343

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 15.537021
$ws.Range("H15").Value = 46.611063
$ws.Range("I15").Value = 0.14592220219851
$ws.Range("J15").Value = 0.14592220219851
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 103.273595
$ws.Range("N15").Value = 309.820785
$ws.Range("O15").Value = 0.5356844639284184
$ws.Range("P15").Value = 0.5356844639284185
$ws.Range("Q15").Value = 1604.564014260495
$ws.Range("R15").Value = 14441.07612834446
$ws.Range("S15").Value = 0.07816825665996312
$ws.Range("T15").Value = 0.07816825665996312

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 15.537021
$ws.Range("H16").Value = 46.611063
$ws.Range("I16").Value = 0.14592220219851
$ws.Range("J16").Value = 0.14592220219851
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 53.963124
$ws.Range("N16").Value = 161.889372
$ws.Range("O16").Value = 0.2799089849815219
$ws.Range("P16").Value = 0.2799089849815219
$ws.Range("Q16").Value = 838.4261908136041
$ws.Range("R16").Value = 7545.835717322437
$ws.Range("S16").Value = 0.04084493550365334
$ws.Range("T16").Value = 0.04084493550365333

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 15.537021
$ws.Range("H17").Value = 46.611063
$ws.Range("I17").Value = 0.14592220219851
$ws.Range("J17").Value = 0.14592220219851
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 33.77521133333333
$ws.Range("N17").Value = 101.325634
$ws.Range("O17").Value = 0.1751934362037625
$ws.Range("P17").Value = 0.1751934362037625
$ws.Range("Q17").Value = 524.766167765438
$ws.Range("R17").Value = 4722.895509888942
$ws.Range("S17").Value = 0.0255646120215772
$ws.Range("T17").Value = 0.0255646120215772
